$d = $word.ActiveDocument

# 1. Sprint Number: 3 -> 2
$d.Content.Find.Execute("Sprint Number: 3", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Sprint Number: 2", 2)

# 2. "Writing unit tests" (Answer for "What did you accomplish yesterday?") ->
#    "Analyze the functions to write unit tests"
$d.Content.Find.Execute("Answer: Writing unit tests", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Answer: Analyze the functions to write unit tests", 2)

# 3. "Continue writing unit tests" (Answer for "What will you do today?") ->
#    "Analyze the functions to write unit tests"
$d.Content.Find.Execute("Continue writing unit tests", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Analyze the functions to write unit tests", 2)
